# Applies the compliance re-check refactor:
#  - "Compliance Overview": violation counts for RDVI drop from 9 to 3
#    (both the per-fund row and the totals row).
#  - "Compliance Details": six checks that used to FAIL now PASS (with
#    violation counts going from 1 to 0), and the red FAIL highlighting
#    (fill style) is cleared back to the workbook's default style, just
#    like the already-passing rows (e.g. row 7).

$wb = $excel.ActiveWorkbook

# --- Sheet: Compliance Overview ---------------------------------------
$overview = $wb.Sheets.Item("Compliance Overview")

# Row 9 (RDVI) : Violations Before / After -> 3
$overview.Range("C9").Value = 3
$overview.Range("D9").Value = 3

# Row 10 (Totals) : Violations Before / After -> 3
$overview.Range("C10").Value = 3
$overview.Range("D10").Value = 3

# --- Sheet: Compliance Details -----------------------------------------
$details = $wb.Sheets.Item("Compliance Details")

# Rows whose check flipped from FAIL to PASS (0 violations), and whose
# red "failed" fill should be cleared to the default (unstyled) look.
$passRows = 2, 5, 6, 8, 9, 10

foreach ($row in $passRows) {
    $details.Range("C$row").Value = "PASS"
    $details.Range("D$row").Value = "PASS"
    $details.Range("E$row").Value = 0
    $details.Range("F$row").Value = 0

    # Clear the red FAIL fill/style on the whole row's data cells so it
    # matches the look of rows that already pass (e.g. row 7).
    $details.Range("A" + $row + ":G" + $row).ClearFormats()
}
